# Daily attendance processing - 2026-01-04 21:56:37
# Swap the order of names in the "Taken By" (column G) cells from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row where that exact value occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
